$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138 (existing rows 138-163 shift down to 139-164),
# then populate the new row with the latest weekly price entry.
$ws.Rows.Item(138).Insert()

$ws.Cells.Item(138, 1).Value = 8
$ws.Cells.Item(138, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(138, 3).Value = "Coquimbo"
$ws.Cells.Item(138, 4).Value = 44798
$ws.Cells.Item(138, 5).Value = 4
$ws.Cells.Item(138, 6).Value = 100112044
$ws.Cells.Item(138, 7).Value = "Perejil"
$ws.Cells.Item(138, 8).Value = "Sin especificar"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 2000
$ws.Cells.Item(138, 11).Value = 2000
$ws.Cells.Item(138, 12).Value = 2500
$ws.Cells.Item(138, 13).Value = 2250
$ws.Cells.Item(138, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(138, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(138, 16).Value = 1500
$ws.Cells.Item(138, 17).Value = 1.5
$ws.Cells.Item(138, 18).Value = "Hortaliza"
